$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new column of data for year 2021 (column R) to the table ---

# R3: blank separator cell above the header row - copy formatting from the
# neighboring Q3 cell (thin/medium bottom border, no value).
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("R3").PasteSpecial(-4122) | Out-Null

# R4: header cell with the year value 2021 - copy formatting from Q4 (2020).
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null
$ws.Range("R4").Value = 2021

# R5: data value for the first indicator - copy formatting from Q5, then
# set the 2021 value and switch the vertical alignment to "top".
$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null
$ws.Range("R5").Value = 0.9
$ws.Range("R5").VerticalAlignment = -4160

# R6: data value for the second indicator - copy formatting from Q6
# (includes the "0.0" number format), then set the 2021 value.
$ws.Range("Q6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null
$ws.Range("R6").Value = 6.5

# Update the sheet's remembered selection to match the saved workbook state.
$ws.Range("T5").Select() | Out-Null
